# Modified controller to pick the one test script.
# Only User_01 (row 2) stays marked to execute; User_02 (row 3) and
# User_03 (row 4) are switched from "Yes" to "No" on the UserModule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserModule")

$ws.Range("D3").Value = "No"
$ws.Range("D4").Value = "No"

$ws.Activate()
$ws.Range("D4").Select()
